# Refactor: Implement services for admin functionality
#
# The underlying change is a small data/UI-state edit inside testImport.xlsx:
#   - the "ImageUrl" value stored in D2 (a shared string) points to a new
#     generated image file name
#   - the workbook's last on-screen selection moves from D8 to D6
#   - (cosmetic, session-only metadata: last-saved path, co-authoring
#     revision id and window position also changed, but those are written by
#     Excel itself from live session/OS state and aren't reachable through
#     the Workbook/Worksheet/Range object model)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ImageUrl cell (D2) to the new generated file name.
$ws.Range("D2").Value = "\Img\Drink\0154118d-a578-4d20-a6ef-963695d6dd7c.jpg"

# Move the active selection from D8 to D6.
$ws.Range("D6").Select()

# Best-effort: nudge the window position to match the recorded session
# (mirrors real Excel's Window.Left/Top, which back the xWindow/yWindow
# values written into <bookViews><workbookView>).
$win = $excel.ActiveWindow
$win.Left = 11025
$win.Top = 11610
